$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new data point for 2026/01/31 (Sat) was missing from the log; insert it
# as a new row 727, pushing the existing rows 727-768 down to 728-769.
$ws.Rows.Item(727).Insert()

# Column A holds plain date-like text (e.g. "2026/01/30"), not a real Excel
# date. Assigning it straight to .Value would get auto-parsed into a date
# serial, so prefix with an apostrophe to force text, exactly like typing
# '2026/01/31 into the cell in the Excel UI.
$ws.Cells.Item(727, 1).Value = "'2026/01/31"
$ws.Cells.Item(727, 2).Value = "土"
$ws.Cells.Item(727, 3).Value = 3
$ws.Cells.Item(727, 4).Value = 26
